# Populated all files with crossref tags, when applicable
# Update the 'Authors' (E) and 'Other found locations' (I) columns for each
# reference row to append the CROSSREF-aware author lists/tags.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Jianbo%Lai%NULL%0,                   Simeng%Ma%NULL%2,                   Ying%Wang%NULL%3,                   Zhongxiang%Cai%NULL%2,                   Jianbo%Hu%NULL%1,                   Ning%Wei%NULL%1,                   Jiang%Wu%NULL%1,                   Hui%Du%NULL%0,                   Tingting%Chen%NULL%1,                   Ruiting%Li%NULL%3,                   Huawei%Tan%NULL%1,                   Lijun%Kang%NULL%2,                   Lihua%Yao%NULL%1,                   Manli%Huang%NULL%1,                   Huafen%Wang%NULL%1,                   Gaohua%Wang%NULL%1,                   Zhongchun%Liu%NULL%2,                   Shaohua%Hu%NULL%1]'
$ws.Range("I2").Value = '_PMC_CROSSREF'
$ws.Range("E3").Value = '[Jiang%Du%NULL%2,                   Lu%Dong%NULL%1,                   Tao%Wang%NULL%0,                   Chenxin%Yuan%NULL%1,                   Rao%Fu%NULL%1,                   Lei%Zhang%NULL%2,                   Bo%Liu%NULL%2,                   Mingmin%Zhang%NULL%1,                   Yuanyuan%Yin%NULL%1,                   Jiawen%Qin%NULL%1,                   Jennifer%Bouey%NULL%1,                   Min%Zhao%NULL%0,                   Xin%Li%NULL%0]'
$ws.Range("I3").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E4").Value = '[Guo%Li%xref no email%0,            Jinfeng%Miao%xref no email%0,            Hui%Wang%xref no email%0,            Shabei%Xu%xref no email%0,            Wenzhe%Sun%xref no email%0,            Yebin%Fan%xref no email%0,            Chenyan%Zhang%xref no email%0,            Suiqiang%Zhu%xref no email%0,            Zhou%Zhu%xref no email%0,            Wei%Wang%xref no email%0]'
$ws.Range("I4").Value = '_CROSSREF'
$ws.Range("E5").Value = '[Chen-Yun%Liu%NULL%2,                   Yun-zhi%Yang%NULL%1,                   Xiao-Ming%Zhang%NULL%1,                   Xinying%Xu%NULL%2,                   Xinying%Xu%NULL%0,                   Qing-Li%Dou%NULL%1,                   Wen-Wu%Zhang%NULL%1,                   Andy S. K.%Cheng%NULL%1]'
$ws.Range("I5").Value = '_PMC_CROSSREF'
$ws.Range("E6").Value = '[Wen%Lu%luwen67@sina.com%0,                   Hang%Wang%NULL%1,                   Yuxing%Lin%NULL%1,                   Li%Li%lilifuzhou@126.com%0]'
$ws.Range("I6").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E7").Value = '[Jing%Qi%NULL%1,                   Jing%Xu%NULL%0,                   Bo-Zhi%Li%NULL%1,                   Jin-Sha%Huang%NULL%1,                   Yuan%Yang%NULL%2,                   Zhen-Tao%Zhang%NULL%1,                   Dong-Ai%Yao%NULL%1,                   Qun-Hui%Liu%NULL%1,                   Min%Jia%NULL%1,                   Dao-Kai%Gong%NULL%1,                   Xiao-Hong%Ni%NULL%1,                   Qi-Mei%Zhang%NULL%1,                   Fu-Rong%Shang%NULL%1,                   Nian%Xiong%NULL%0,                   Chun-Li%Zhu%NULL%1,                   Tao%Wang%wangtaowh@hust.edu.cn%0,                   Xi%Zhang%SMRC301@163.com%1]'
$ws.Range("I7").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E9").Value = '[Xingyue%Song%NULL%0,                   Wenning%Fu%NULL%1,                   Xiaoran%Liu%NULL%1,                   Zhiqian%Luo%NULL%1,                   Rixing%Wang%NULL%1,                   Ning%Zhou%NULL%0,                   Shijiao%Yan%NULL%1,                   Chuanzhu%Lv%NULL%1]'
$ws.Range("I9").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E10").Value = '[S%Wang%NULL%1,                   L%Xie%NULL%1,                   Y%Xu%NULL%1,                   S%Yu%NULL%1,                   B%Yao%professoryao@aliyun.com%1,                   D%Xiang%xiangdannuli@163.com%1]'
$ws.Range("I10").Value = '_PMC_CROSSREF'
$ws.Range("E11").Value = '[Yeen%Huang%NULL%2,                   Ning%Zhao%zhaoning2018@email.szu.edu.cn%1]'
$ws.Range("I11").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E13").Value = '[Wen-rui%Zhang%NULL%0,                   Kun%Wang%NULL%1,                   Lu%Yin%NULL%1,                   Wen-feng%Zhao%NULL%1,                   Qing%Xue%NULL%1,                   Mao%Peng%NULL%1,                   Bao-quan%Min%NULL%1,                   Qing%Tian%NULL%1,                   Hai-xia%Leng%NULL%1,                   Jia-lin%Du%NULL%1,                   Hong%Chang%NULL%1,                   Yuan%Yang%NULL%0,                   Wei%Li%NULL%4,                   Fang-fang%Shangguan%NULL%1,                   Tian-yi%Yan%NULL%1,                   Hui-qing%Dong%NULL%1,                   Ying%Han%NULL%1,                   Yu-ping%Wang%NULL%1,                   Fiammetta%Cosci%NULL%1,                   Hong-xing%Wang%NULL%1]'
$ws.Range("I13").Value = '_PMC_CROSSREF'
$ws.Range("E14").Value = '[Chenxi%Zhang%NULL%1,                   Lulu%Yang%NULL%1,                   Shuai%Liu%NULL%1,                   Simeng%Ma%NULL%0,                   Ying%Wang%NULL%0,                   Zhongxiang%Cai%NULL%0,                   Hui%Du%NULL%0,                   Ruiting%Li%NULL%0,                   Lijun%Kang%NULL%0,                   Meilei%Su%NULL%1,                   Jihui%Zhang%NULL%1,                   Zhongchun%Liu%NULL%0,                   Bin%Zhang%NULL%0]'
$ws.Range("I14").Value = '_PMC_CROSSREF'
$ws.Range("E15").Value = '[Xiao%Xiao%NULL%2,                   Xiaobin%Zhu%NULL%1,                   Shuai%Fu%NULL%1,                   Yugang%Hu%NULL%1,                   Xiaoning%Li%NULL%1,                   Jinsong%Xiao%NULL%1]'
$ws.Range("I15").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E16").Value = '[Yun%Chen%NULL%2,                   Hao%Zhou%haoye320@163.com%1,                   Yan%Zhou%NULL%1,                   Fang%Zhou%NULL%1]'
$ws.Range("I16").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E17").Value = '[Zhi-hao%Tu%NULL%1,                   Jing-wen%He%NULL%2,                   Jing-wen%He%NULL%0,                   Na%Zhou%NULL%1,                   Oliver%Schildgen.%NULL%3,                   Oliver%Schildgen.%NULL%0,                   Oliver%Schildgen.%NULL%0]'
$ws.Range("I17").Value = '_PMC_CROSSREF'
$ws.Range("E18").Value = '[Juhong%Zhu%NULL%2,                   Lin%Sun%NULL%1,                   Lan%Zhang%NULL%1,                   Huan%Wang%NULL%1,                   Ajiao%Fan%NULL%1,                   Bin%Yang%NULL%1,                   Wei%Li%NULL%0,                   Shifu%Xiao%NULL%1]'
$ws.Range("I18").Value = '_PMC_CROSSREF'
$ws.Range("E19").Value = '[Seshadri Sekhar%Chatterjee%NULL%2,                   Ranjan%Bhattacharyya%NULL%1,                   Sumita%Bhattacharyya%NULL%1,                   Sukanya%Gupta%NULL%1,                   Soumitra%Das%NULL%1,                   Bejoy Bikram%Banerjee%NULL%1]'
$ws.Range("I19").Value = '_PMC_CROSSREF'
$ws.Range("E20").Value = '[Simmi%Gupta%NULL%1,                   Kinshuk%Kohli%NULL%1,                   P%Padmakumari%NULL%1,                   Prashant Kumar%Dixit%NULL%1,                   Ajay Shankar%Prasad%NULL%1,                   Bharath Seshadri%Chakravarthy%NULL%1,                   Rajat%Shukla%NULL%1,                   Parthasarathi%Ghana%NULL%1,                   Debasish%Mahapatra%NULL%1,                   G%Varadaraj%NULL%1]'
$ws.Range("I20").Value = '_PMC_CROSSREF'
$ws.Range("E21").Value = '[William%Wilson%NULL%1,                   Jeffrey Pradeep%Raj%NULL%1,                   Seema%Rao%NULL%1,                   Murtuza%Ghiya%NULL%1,                   Nisanth Menon%Nedungalaparambil%NULL%1,                   Harshit%Mundra%NULL%1,                   Roshan%Mathew%NULL%1]'
$ws.Range("I21").Value = '_PMC_CROSSREF'
$ws.Range("E22").Value = '[Alexander%Muacevic%NULL%0,                   John R%Adler%NULL%0,                   Ram%Sandesh%NULL%2,                   Ram%Sandesh%NULL%0,                   Wajeeha%Shahid%NULL%1,                   Kapeel%Dev%NULL%1,                   Nikeeta%Mandhan%NULL%1,                   Prem%Shankar%NULL%1,                   Anam%Shaikh%NULL%1,                   Amber%Rizwan%NULL%1]'
$ws.Range("I22").Value = '_PMC_CROSSREF'
$ws.Range("E24").Value = '[Rümeysa Yeni%Elbay%rumeysa.yenielbay@medeniyet.edu.tr%2,                   Ayşe%Kurtulmuş%NULL%1,                   Selim%Arpacıoğlu%NULL%1,                   Emrah%Karadere%NULL%1]'
$ws.Range("I24").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E25").Value = '[Benjamin Y.Q.%Tan%NULL%0,                   Nicholas W.S.%Chew%NULL%1,                   Grace K.H.%Lee%NULL%2,                   Mingxue%Jing%NULL%2,                   Yihui%Goh%NULL%2,                   Leonard L.L.%Yeo%NULL%2,                   Ka%Zhang%NULL%1,                   Howe-Keat%Chin%NULL%1,                   Aftab%Ahmad%NULL%2,                   Faheem Ahmed%Khan%NULL%1,                   Ganesh Napolean%Shanmugam%NULL%1,                   Bernard P.L.%Chan%NULL%2,                   Sibi%Sunny%NULL%2,                   Bharatendu%Chandra%NULL%2,                   Jonathan J.Y.%Ong%NULL%2,                   Prakash R.%Paliwal%NULL%2,                   Lily Y.H.%Wong%NULL%2,                   Renarebecca%Sagayanathan%NULL%2,                   Jin Tao%Chen%NULL%2,                   Alison Ying%Ying Ng%NULL%2,                   Hock Luen%Teoh%NULL%2,                   Cyrus S.%Ho%NULL%2,                   Roger C.%Ho%NULL%2,                   Vijay K.%Sharma%NULL%2]'
$ws.Range("I25").Value = '_PMC_CROSSREF'
$ws.Range("E26").Value = '[Alessandra%Rodrigues Martins%coreGivesNoEmail%1,                  Aline%Cristina Correa Anselmo%coreGivesNoEmail%1,                  Assun\u00e7\u00e3o%AA%coreGivesNoEmail%1,                  Botega%NJ%coreGivesNoEmail%1,                  Eduardo%Bassani Dal\u2019Bosco%coreGivesNoEmail%1,                  Freitas%FMB%coreGivesNoEmail%1,                  Guilherme%Arcaro%coreGivesNoEmail%1,                  Kang%L%coreGivesNoEmail%1,                  Lai%J%coreGivesNoEmail%1,                  Lana%RM%coreGivesNoEmail%1,                  Lara%Simone Messias Floriano%coreGivesNoEmail%1,                  Leonelli%LB%coreGivesNoEmail%1,                  Le\u00e3o%AM%coreGivesNoEmail%1,                  Liu%Y%coreGivesNoEmail%1,                  Martins%C%coreGivesNoEmail%1,                  Moura%A%coreGivesNoEmail%1,                  Oliveira%FP%coreGivesNoEmail%1,                  Pappa%S%coreGivesNoEmail%1,                  Rodriguez%EOL%coreGivesNoEmail%1,                  Ruback%SP%coreGivesNoEmail%1,                  Sena%AFJ%coreGivesNoEmail%1,                  Silva%DSD%coreGivesNoEmail%2,                  Silva%DSD%coreGivesNoEmail%0,                  Spoorthy%MS%coreGivesNoEmail%1,                  Suellen%Vienscoski Skupien%coreGivesNoEmail%1,                  Trettene%AS%coreGivesNoEmail%1,                  Vieira%NF%coreGivesNoEmail%1,                  Wang%SM%coreGivesNoEmail%1,                  Zigmond%AS%coreGivesNoEmail%1]'
$ws.Range("I26").Value = '_CROSSREF'
$ws.Range("E27").Value = '[Emanuele Maria%Giusti%NULL%0,                   Elisa%Pedroli%NULL%1,                   Guido E.%D''Aniello%NULL%1,                   Chiara%Stramba Badiale%NULL%1,                   Giada%Pietrabissa%NULL%1,                   Chiara%Manna%NULL%1,                   Marco%Stramba Badiale%NULL%1,                   Giuseppe%Riva%NULL%1,                   Gianluca%Castelnuovo%NULL%1,                   Enrico%Molinari%NULL%1]'
$ws.Range("I27").Value = '_PMC_CROSSREF'
$ws.Range("E28").Value = '[Jacek C.%Szepietowski%jacek.szepietowski@umed.wroc.pl%1,                   Piotr%Krajewski%NULL%2,                   Piotr%Krajewski%NULL%0,                   Rafał%Biłynicki‐Birula%NULL%2,                   Rafał%Biłynicki‐Birula%NULL%0,                   Paweł%Poznański%NULL%2,                   Paweł%Poznański%NULL%0,                   Magdalena%Krajewska%NULL%2,                   Magdalena%Krajewska%NULL%0,                   Joanna%Rymaszewska%NULL%2,                   Joanna%Rymaszewska%NULL%0,                   Łukasz%Matusiak%NULL%2,                   Łukasz%Matusiak%NULL%0]'
$ws.Range("I28").Value = '_PMC_CROSSREF'
$ws.Range("E29").Value = '[Stephen X.%Zhang%NULL%2,                   Jing%Liu%NULL%0,                   Asghar%Afshar Jahanshahi%NULL%1,                   Khaled%Nawaser%NULL%1,                   Ali%Yousefi%NULL%1,                   Jizhen%Li%NULL%1,                   Shuhua%Sun%NULL%1]'
$ws.Range("I29").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E30").Value = '[Abdallah Y.%Naser%abdallah.naser@iu.edu.jo%2,                   Eman Zmaily%Dahmash%NULL%2,                   Eman Zmaily%Dahmash%NULL%0,                   Rabaa%Al‐Rousan%NULL%2,                   Rabaa%Al‐Rousan%NULL%0,                   Hassan%Alwafi%NULL%1,                   Hamzeh Mohammad%Alrawashdeh%NULL%1,                   Imene%Ghoul%NULL%1,                   Anwer%Abidine%NULL%1,                   Mohammed A.%Bokhary%NULL%1,                   Hadeel T.%AL‐Hadithi%NULL%1,                   Dalia%Ali%NULL%1,                   Rasha%Abuthawabeh%NULL%1,                   Ghada Mohammad%Abdelwahab%NULL%1,                   Yosra J.%Alhartani%NULL%1,                   Haneen%Al Muhaisen%NULL%1,                   Ayah%Dagash%NULL%1,                   Hamad S.%Alyami%NULL%1]'
$ws.Range("I30").Value = '_PMC_CROSSREF'
$ws.Range("E31").Value = '[Anoop Krishna%Gupta%NULL%1,                   Aseem%Mehra%NULL%2,                   Aseem%Mehra%NULL%0,                   Abhash%Niraula%NULL%2,                   Abhash%Niraula%NULL%0,                   Khagendra%Kafle%NULL%2,                   Khagendra%Kafle%NULL%0,                   Saroj Prasad%Deo%NULL%2,                   Saroj Prasad%Deo%NULL%0,                   Babita%Singh%NULL%2,                   Babita%Singh%NULL%0,                   Swapnajeet%Sahoo%NULL%2,                   Swapnajeet%Sahoo%NULL%0,                   Sandeep%Grover%NULL%1]'
$ws.Range("I31").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E32").Value = '[Thomas G.%Kannampallil%NULL%1,                   Charles W.%Goss%NULL%2,                   Charles W.%Goss%NULL%0,                   Bradley A.%Evanoff%NULL%1,                   Jaime R.%Strickland%NULL%1,                   Rebecca P.%McAlister%NULL%1,                   Jennifer%Duncan%NULL%1,                   Michio%Murakami%NULL%2,                   Michio%Murakami%NULL%0]'
$ws.Range("I32").Value = '_PMC_CROSSREF'
$ws.Range("E33").Value = '[Nicholas W.S.%Chew%NULL%1,                   Grace K.H.%Lee%NULL%0,                   Benjamin Y.Q.%Tan%NULL%1,                   Mingxue%Jing%NULL%0,                   Yihui%Goh%NULL%0,                   Nicholas J.H.%Ngiam%NULL%1,                   Leonard L.L.%Yeo%NULL%0,                   Aftab%Ahmad%NULL%0,                   Faheem%Ahmed Khan%NULL%1,                   Ganesh%Napolean Shanmugam%NULL%1,                   Arvind K.%Sharma%NULL%1,                   R.N.%Komalkumar%NULL%1,                   P.V.%Meenakshi%NULL%1,                   Kenam%Shah%NULL%1,                   Bhargesh%Patel%NULL%1,                   Bernard P.L.%Chan%NULL%0,                   Sibi%Sunny%NULL%0,                   Bharatendu%Chandra%NULL%0,                   Jonathan J.Y.%Ong%NULL%0,                   Prakash R.%Paliwal%NULL%0,                   Lily Y.H.%Wong%NULL%0,                   Renarebecca%Sagayanathan%NULL%0,                   Jin Tao%Chen%NULL%0,                   Alison Ying%Ying Ng%NULL%0,                   Hock Luen%Teoh%NULL%0,                   Georgios%Tsivgoulis%NULL%1,                   Cyrus S.%Ho%NULL%0,                   Roger C.%Ho%NULL%0,                   Vijay K.%Sharma%NULL%0]'
$ws.Range("I33").Value = '_PMC_elsevier_CROSSREF'

Write-Output "Updated Authors/Other-found-locations columns for 29 rows"
